# Add 2022-Q4 data
# 1) Insert a new "2022-Q4" worksheet right after the "总计" (summary) sheet.
# 2) Insert a new row into the "总计" sheet for the 2022-Q4 totals, pushing the
#    existing Q3/Q2/Q1 rows down by one and renumbering the index column.
# 3) Populate the new "2022-Q4" sheet with the per-fund holding data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new sheet right after "总计" and name it "2022-Q4".
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Grab the old "2022-Q3" sheet (by name, since indices shift after Add) to use
# as a formatting template for the header row + index column styling.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("A1:H3").Copy($q4.Range("A1:H3"))

# Stamp the data-row formatting (only column A is styled) down through row 10.
for ($r = 4; $r -le 10; $r++) {
    $q4.Range("A2:H2").Copy($q4.Range("A" + $r + ":H" + $r))
}

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2:2").Insert()

# New row 2 needs the same look as the (now shifted) data rows below it
# (Insert() otherwise carries the header row's bold styling down into B2:D2).
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("B3:D3").Copy($summary.Range("B2:D2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 0.67

# Renumber the index column for the rows that got pushed down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Step 3: fill in the 2022-Q4 per-fund data.
# ---------------------------------------------------------------------------
# Re-resolve by name defensively (sheet handles in this engine behave
# positionally, and no further sheet insertions happen after this point, but
# re-fetching by name keeps this block robust regardless).
$q4 = $wb.Worksheets.Item("2022-Q4")

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'010405"
$q4.Range("C2").Value = "惠升医药健康6个月持有期混合"
$q4.Range("D2").Value = "'9.98"
$q4.Range("E2").Value = "'72.67"
$q4.Range("F2").Value = "'2.52"
$q4.Range("G2").Value = "'0.2515"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'011738"
$q4.Range("C3").Value = "华安兴安优选一年持有期混合A"
$q4.Range("D3").Value = "'13.05"
$q4.Range("E3").Value = "'53.71"
$q4.Range("F3").Value = "'1.53"
$q4.Range("G3").Value = "'0.1997"
$q4.Range("H3").Value = 10

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'011739"
$q4.Range("C4").Value = "华安兴安优选一年持有期混合C"
$q4.Range("D4").Value = "'8.07"
$q4.Range("E4").Value = "'53.71"
$q4.Range("F4").Value = "'1.53"
$q4.Range("G4").Value = "'0.1235"
$q4.Range("H4").Value = 10

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'011390"
$q4.Range("C5").Value = "华安添祥6个月持有期混合A"
$q4.Range("D5").Value = "'6.21"
$q4.Range("E5").Value = "'32.32"
$q4.Range("F5").Value = "'1.08"
$q4.Range("G5").Value = "'0.0671"
$q4.Range("H5").Value = 8

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'013920"
$q4.Range("C6").Value = "兴华创新医疗6个月持有混合A"
$q4.Range("D6").Value = "'0.22"
$q4.Range("E6").Value = "'92.79"
$q4.Range("F6").Value = "'6.32"
$q4.Range("G6").Value = "'0.0139"
$q4.Range("H6").Value = 1

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'013921"
$q4.Range("C7").Value = "兴华创新医疗6个月持有混合C"
$q4.Range("D7").Value = "'0.06"
$q4.Range("E7").Value = "'92.79"
$q4.Range("F7").Value = "'6.32"
$q4.Range("G7").Value = "'0.0038"
$q4.Range("H7").Value = 1

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'012315"
$q4.Range("C8").Value = "创金合信港股通成长股票A"
$q4.Range("D8").Value = "'0.12"
$q4.Range("E8").Value = "'89.18"
$q4.Range("F8").Value = "'2.78"
$q4.Range("G8").Value = "'0.0033"
$q4.Range("H8").Value = 10

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "'012316"
$q4.Range("C9").Value = "创金合信港股通成长股票C"
$q4.Range("D9").Value = "'0.11"
$q4.Range("E9").Value = "'89.18"
$q4.Range("F9").Value = "'2.78"
$q4.Range("G9").Value = "'0.0031"
$q4.Range("H9").Value = 10

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "'016181"
$q4.Range("C10").Value = "华安添祥6个月持有期混合C"
$q4.Range("D10").Value = "'0.00"
$q4.Range("E10").Value = "'32.32"
$q4.Range("F10").Value = "'1.08"
$q4.Range("G10").Value = 0
$q4.Range("H10").Value = 8

Write-Host "2022-Q4 sheet added and summary sheet updated"
